# Update cryptos worksheet with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.384.82"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.693.30"
$ws.Range("E4").Value = "  +0.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.20"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5495"
$ws.Range("E6").Value = "  +4.28%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2739"
$ws.Range("E8").Value = "  +1.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06459"
$ws.Range("E9").Value = "  +0.35%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.99"
$ws.Range("E10").Value = "  -0.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07671"
$ws.Range("E11").Value = "  +2.64%  "
$ws.Range("D12").Value = "1.706.84"
$ws.Range("E12").Value = "  +0.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.538"
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5830"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008348"
$ws.Range("E15").Value = "  -2.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.45"
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").Value = "26.424.39"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.938"
$ws.Range("E18").Value = "  -0.48%  "
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("E20").Value = "  +0.58%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "192.07"
$ws.Range("E21").Value = "  +1.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.245"
$ws.Range("E22").Value = "  +0.30%  "
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.98"
$ws.Range("E24").Value = "  +2.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1324"
$ws.Range("E25").Value = "  +7.27%  "
$ws.Range("E26").Value = "  +2.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06290"
$ws.Range("E28").Value = "  -5.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.384"
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.332"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.606"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.042"
$ws.Range("E34").Value = "  +1.14%  "
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.414"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.707"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.204"
$ws.Range("E38").Value = "  -2.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01642"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("D40").Value = "1.115.78"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8888"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.89"
$ws.Range("E43").Value = "  +1.02%  "
$ws.Range("D44").Value = "1.844.70"
$ws.Range("E44").Value = "  +0.32%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.49"
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000109"
$ws.Range("E46").Value = "  -3.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.014"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.209"
$ws.Range("E48").Value = "  +0.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05287"
$ws.Range("E49").Value = "  +0.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4304"
$ws.Range("E50").Value = "  +0.04%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.099"
$ws.Range("E51").Value = "  +0.49%  "

Write-Output "Updated $($wb.ActiveSheet.Name) with latest crypto data"
